# SARAALERT-1260: Allow vaccine table to be populated on import (#958)
# Adds two "Vaccine N ..." blocks of columns (CO:CX) to the Monitorees
# import-format fixture sheet, with sample data for rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row (row 1): new column headers CO1:CX1
# ---------------------------------------------------------------------
$ws.Range("CO1").Value = "Vaccine 1 Group Name"
$ws.Range("CP1").Value = "Vaccine 1 Product Name"
$ws.Range("CQ1").NumberFormat = "@"
$ws.Range("CQ1").Value = "Vaccine 1 Administration Date"
$ws.Range("CR1").Value = "Vaccine 1 Dose Number"
$ws.Range("CS1").Value = "Vaccine 1 Notes"
$ws.Range("CT1").Value = "Vaccine 2 Group Name"
$ws.Range("CU1").Value = "Vaccine 2 Product Name"
$ws.Range("CV1").NumberFormat = "@"
$ws.Range("CV1").Value = "Vaccine 2 Administration Date"
$ws.Range("CW1").Value = "Vaccine 2 Dose Number"
$ws.Range("CX1").Value = "Vaccine 2 Notes"

# ---------------------------------------------------------------------
# Row 2
# ---------------------------------------------------------------------
$ws.Range("CO2").Value = "COVID-19"
$ws.Range("CP2").Value = "Moderna COVID-19 Vaccine"
$ws.Range("CQ2").NumberFormat = "@"
$ws.Range("CQ2").Value = "2020-06-01"
$ws.Range("CR2").Value = 1
$ws.Range("CS2").Value = "notes 1"
$ws.Range("CT2").Value = "COVID-19"
$ws.Range("CU2").Value = "Moderna COVID-19 Vaccine"
$ws.Range("CV2").NumberFormat = "@"
$ws.Range("CV2").Value = "2020-06-20"
$ws.Range("CW2").Value = 2
$ws.Range("CX2").Value = "notes 2"

# ---------------------------------------------------------------------
# Row 3
# ---------------------------------------------------------------------
$ws.Range("CO3").Value = "COVID-19"
$ws.Range("CP3").Value = "Pfizer-BioNTech COVID-19 Vaccine"
$ws.Range("CQ3").NumberFormat = "@"
$ws.Range("CQ3").Value = "2020-06-02"
$ws.Range("CR3").Value = 1
$ws.Range("CT3").Value = "COVID-19"
$ws.Range("CU3").Value = "Pfizer-BioNTech COVID-19 Vaccine"
$ws.Range("CV3").NumberFormat = "@"
$ws.Range("CV3").Value = "2020-06-21"
$ws.Range("CW3").Value = 2

# ---------------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------------
$ws.Range("CO4").Value = "COVID-19"
$ws.Range("CP4").Value = "Unknown"
$ws.Range("CQ4").NumberFormat = "@"
$ws.Range("CQ4").Value = "2020-06-04"
$ws.Range("CR4").Value = 1
$ws.Range("CT4").Value = "COVID-19"
$ws.Range("CU4").Value = "Unknown"
$ws.Range("CV4").NumberFormat = "@"
$ws.Range("CV4").Value = "2020-06-22"
$ws.Range("CW4").Value = 2

# ---------------------------------------------------------------------
# Row 5
# ---------------------------------------------------------------------
$ws.Range("CO5").Value = "COVID-19"
$ws.Range("CP5").Value = "Moderna COVID-19 Vaccine"
$ws.Range("CQ5").NumberFormat = "@"
$ws.Range("CQ5").Value = "2020-06-01"
$ws.Range("CR5").Value = 1

# ---------------------------------------------------------------------
# Row 6 - CO6 previously held an empty, date-styled placeholder cell;
# clear it first so the stale style is dropped before the new text goes in.
# ---------------------------------------------------------------------
$ws.Range("CO6").Clear()
$ws.Range("CO6").Value = "COVID-19"
$ws.Range("CP6").Value = "Janssen (J&J) COVID-19 Vaccine"
$ws.Range("CQ6").NumberFormat = "@"
$ws.Range("CQ6").Value = "2020-06-03"
$ws.Range("CR6").Value = 1

# ---------------------------------------------------------------------
# Row 7 - CR7/CS7 previously held empty, date-styled placeholder cells.
# ---------------------------------------------------------------------
$ws.Range("CS7").Clear()
$ws.Range("CR7").Clear()
$ws.Range("CO7").Value = "COVID-19"
$ws.Range("CP7").Value = "Unknown"
$ws.Range("CQ7").NumberFormat = "@"
$ws.Range("CQ7").Value = "2020-06-02"
$ws.Range("CR7").Value = 1

# ---------------------------------------------------------------------
# Row 10 - CO10/CR10/CS10 previously held empty, date-styled placeholder
# cells which are simply removed (no vaccine data for this row).
# ---------------------------------------------------------------------
$ws.Range("CO10").Clear()
$ws.Range("CR10").Clear()
$ws.Range("CS10").Clear()

# ---------------------------------------------------------------------
# Column widths for the new vaccine columns (best match of the widths
# Excel computed for this content).
# ---------------------------------------------------------------------
$ws.Columns("CO").ColumnWidth = 20.33203125
$ws.Columns("CP").ColumnWidth = 31
$ws.Columns("CQ").ColumnWidth = 25.6640625
$ws.Columns("CR").ColumnWidth = 21.1640625
$ws.Columns("CS").ColumnWidth = 14.5
$ws.Columns("CT").ColumnWidth = 20.33203125
$ws.Columns("CU").ColumnWidth = 31
$ws.Columns("CV").ColumnWidth = 25.6640625
$ws.Columns("CW").ColumnWidth = 21.1640625
$ws.Columns("CX").ColumnWidth = 14.5

# ---------------------------------------------------------------------
# Reset the view back to the top-left/default selection (matches the
# saved file no longer being scrolled to CG1 with CQ15 selected).
# ---------------------------------------------------------------------
$ws.Range("A1").Select()
